# feat: add 2022-Q4 data
#
# Before:  sheet "总计" (totals) + sheet "2022-Q3" (fund snapshot for Q3)
# After:   sheet "总计" gains a new top row for 2022-Q4 (old Q3 row pushed
#          down to row 3), the existing "2022-Q3" sheet is turned into
#          "2022-Q4" with refreshed figures, and a fresh "2022-Q3" sheet
#          (an exact copy of the original Q3 snapshot) is inserted right
#          after it.

$wb = $excel.ActiveWorkbook

# Helper: write $value into $range as a genuine TEXT cell (no leading
# apostrophe, no numeric auto-coercion), leaving the cell's style exactly
# as it was before the call (the "@" number-format trick would otherwise
# stamp a stray numFmt onto the cell, so we paste the original formatting
# straight back on top from $styleSource once the text is in place).
function Set-TextValue {
    param($range, $value, $styleSource)
    $range.NumberFormat = "@"
    $range.Value = $value
    $styleSource.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the old 2022-Q3 total down to row 3 and put the
#    new 2022-Q4 total in row 2.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# New row 3 = what row 2 used to hold (2022-Q3 / 1 / 0.13)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.13
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)  # xlPasteFormats -> picks up style "s=2"

# Row 2 becomes the 2022-Q4 total
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.14

# ---------------------------------------------------------------------
# 2) Duplicate the original "2022-Q3" sheet (placed right after it) so the
#    untouched snapshot survives under its original name/format, then
#    rename the original to "2022-Q4" and refresh its numbers.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsQ3)            # -> new sheet "2022-Q3 (2)" right after it
$wsNewQ3 = $wb.Worksheets.Item("2022-Q3 (2)")

$wsQ3.Name = "2022-Q4"
$wsNewQ3.Name = "2022-Q3"

# Re-style the header row + A2 on the "2022-Q4" sheet to style "s=2"
# (matches the "总计" sheet's header look) by pasting format from a cell
# that already carries it.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2").PasteSpecial(-4122)

# Updated Q4 figures (fund code / name stay the same - only the numbers
# that come from the quarterly snapshot move)
$blankQ4 = $wsQ3.Range("Z100")   # untouched cell -> plain/default style "0"
Set-TextValue $wsQ3.Range("D2") "15.29"  $blankQ4
Set-TextValue $wsQ3.Range("E2") "76.11"  $blankQ4
Set-TextValue $wsQ3.Range("F2") "0.93"   $blankQ4
Set-TextValue $wsQ3.Range("G2") "0.1422" $blankQ4
$wsQ3.Range("H2").Value = 8

Write-Host "2022-Q4 / 2022-Q3 split complete"
